$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" ---
$wsMeta = $wb.Worksheets.Item("Metadata")

# URL: ibm.com -> linuxforhealth.org
$wsMeta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/sexual-orientation"

# Version: 7.0.0 -> 8.0.0
$wsMeta.Range("B3").Value = "8.0.0"

# Date: 2022-09-08T16:11:15+00:00 -> 2022-11-10T16:00:46+00:00
$wsMeta.Range("B8").Value = "2022-11-10T16:00:46+00:00"

# Publisher: Alvearie Team -> LinuxForHealth Team
$wsMeta.Range("B9").Value = "LinuxForHealth Team"

# --- Sheet "Include from Sexual Orientati" ---
$wsInclude = $wb.Worksheets.Item("Include from Sexual Orientati")

# System URI: ibm.com -> linuxforhealth.org
$wsInclude.Range("B4").Value = "http://linuxforhealth.org/fhir/cdm/CodeSystem/sexual-orientation"
